$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultat_etape")

# ADJ row: clarify the label for "ADJ" from the generic "Admis" to "Admis par Jury"
$ws.Range("C8").Value = "Admis par Jury"

# ADMP row: clarify the label for "ADMP" from the generic "Admis" to "ADM à poursuiv"
$ws.Range("C11").Value = "ADM à poursuiv"
